# "Atualizacao de bases das ligas" - swap the match-row data between the
# two rows in each pair below. The leading sequential index in column A
# (37/38 and 109/110) stays put; everything from column B (id) through
# column AC (PL_AhUnder) trades places between the two rows of the pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($RowA, $RowB) {
    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")

    $dataA = $rangeA.Value2
    $dataB = $rangeB.Value2

    $rangeA.Value2 = $dataB
    $rangeB.Value2 = $dataA
}

Swap-RowData 39 40
Swap-RowData 111 112
